$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '69.967.59'
$ws.Cells.Item(2, 5).Value = '  -0.04%  '

$ws.Cells.Item(3, 4).Value = '3.541.25'
$ws.Cells.Item(3, 5).Value = '  +0.02%  '

$ws.Cells.Item(4, 5).Value = '  -0.20%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '603.67'
$ws.Cells.Item(5, 5).Value = '  -1.95%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '196.50'
$ws.Cells.Item(6, 5).Value = '  +4.36%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.634'
$ws.Cells.Item(7, 5).Value = '  -0.31%  '

$ws.Cells.Item(9, 5).Value = '  -4.31%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.653'
$ws.Cells.Item(10, 5).Value = '  -1.71%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '53.95'
$ws.Cells.Item(11, 5).Value = '  +0.17%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000303'
$ws.Cells.Item(12, 5).Value = '  -1.79%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '9.55'
$ws.Cells.Item(13, 5).Value = '  -2.44%  '

$ws.Cells.Item(14, 4).Value = '4.096.03'
$ws.Cells.Item(14, 5).Value = '  -0.34%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '607.23'
$ws.Cells.Item(15, 5).Value = '  -1.31%  '

$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '19.21'
$ws.Cells.Item(16, 5).Value = '  -0.06%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '70.108.36'
$ws.Cells.Item(17, 5).Value = '  +0.07%  '

$ws.Cells.Item(18, 2).Value = 'Uniswap'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.78'
$ws.Cells.Item(18, 5).Value = '  -0.97%  '

$ws.Cells.Item(19, 4).Value = '3.533.86'
$ws.Cells.Item(19, 5).Value = '  -0.72%  '

$ws.Cells.Item(20, 5).Value = '  +0.41%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.997'
$ws.Cells.Item(21, 5).Value = '  -0.50%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '18.15'
$ws.Cells.Item(22, 5).Value = '  +3.19%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.32'
$ws.Cells.Item(23, 5).Value = '  +5.02%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '102.83'
$ws.Cells.Item(24, 5).Value = '  -2.66%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '4.62'
$ws.Cells.Item(25, 5).Value = '  -1.67%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.16'
$ws.Cells.Item(26, 5).Value = '  +3.91%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.92'
$ws.Cells.Item(27, 5).Value = '  -0.49%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.65'
$ws.Cells.Item(28, 5).Value = '  -4.37%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '33.71'
$ws.Cells.Item(29, 5).Value = '  -1.55%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.15'
$ws.Cells.Item(30, 5).Value = '  +0.70%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.34'
$ws.Cells.Item(31, 5).Value = '  +16.09%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '12.65'
$ws.Cells.Item(32, 5).Value = '  +0.93%  '

$ws.Cells.Item(33, 5).Value = '  -1.45%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '63.22'
$ws.Cells.Item(34, 5).Value = '  -1.58%  '

$ws.Cells.Item(35, 4).Value = '0.0₃0862'
$ws.Cells.Item(35, 5).Value = '  +10.34%  '

$ws.Cells.Item(36, 4).Value = '3.735.29'
$ws.Cells.Item(36, 5).Value = '  +5.22%  '

$ws.Cells.Item(37, 5).Value = '  +0.01%  '

$ws.Cells.Item(38, 2).Value = 'Fetch.AI'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.04'
$ws.Cells.Item(38, 5).Value = '  -3.66%  '

$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.64'
$ws.Cells.Item(39, 5).Value = '  +0.77%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.394'
$ws.Cells.Item(40, 5).Value = '  -1.46%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '36.73'
$ws.Cells.Item(41, 5).Value = '  -1.04%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '489.27'
$ws.Cells.Item(42, 5).Value = '  -9.17%  '

$ws.Cells.Item(43, 5).Value = '  -4.91%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0458'
$ws.Cells.Item(44, 5).Value = '  -1.39%  '

$ws.Cells.Item(45, 5).Value = '  -1.94%  '

$ws.Cells.Item(46, 5).Value = '  -4.26%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.31'
$ws.Cells.Item(47, 5).Value = '  -1.20%  '

$ws.Cells.Item(48, 5).Value = '  +0.29%  '

$ws.Cells.Item(49, 5).Value = '  -4.02%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.000256'
$ws.Cells.Item(50, 5).Value = '  +5.56%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '130.96'
$ws.Cells.Item(51, 5).Value = '  -1.18%  '
